# Swap the contents of row 25 and row 26 (only the columns that actually
# differ between the two records) in the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $addr25 = "$col" + "25"
    $addr26 = "$col" + "26"

    $val25 = $ws.Range($addr25).Value2
    $val26 = $ws.Range($addr26).Value2

    $ws.Range($addr25).Value2 = $val26
    $ws.Range($addr26).Value2 = $val25
}
